$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to Text format while writing values so that
# numeric-looking strings (e.g. "301.17", "0.08138", "150.00") are not
# auto-converted to numbers by Excel, then restore the original style.
$origStyleD = $ws.Range("D2:D51").Style
$origStyleE = $ws.Range("E2:E51").Style
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "23.127.91"
$ws.Range("E2").Value = "  -3.33%  "

$ws.Range("D3").Value = "1.600.74"
$ws.Range("E3").Value = "  -2.79%  "

$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("D6").Value = "301.17"
$ws.Range("E6").Value = "  -2.84%  "

$ws.Range("D7").Value = "0.3783"
$ws.Range("E7").Value = "  -2.61%  "

$ws.Range("D8").Value = "0.3644"
$ws.Range("E8").Value = "  -3.96%  "

$ws.Range("D9").Value = "50.10"
$ws.Range("E9").Value = "  -3.62%  "

$ws.Range("E10").Value = "  -6.19%  "

$ws.Range("E11").Value = "  -0.25%  "

$ws.Range("D12").Value = "0.08138"
$ws.Range("E12").Value = "  -3.46%  "

$ws.Range("D13").Value = "22.97"
$ws.Range("E13").Value = "  -3.58%  "

$ws.Range("D14").Value = "6.579"
$ws.Range("E14").Value = "  -6.45%  "

$ws.Range("E15").Value = "  -3.55%  "

$ws.Range("D16").Value = "7.363"
$ws.Range("E16").Value = "  -7.63%  "

$ws.Range("D17").Value = "1.602.39"
$ws.Range("E17").Value = "  -2.90%  "

$ws.Range("D18").Value = "91.70"
$ws.Range("E18").Value = "  -2.56%  "

$ws.Range("D19").Value = "0.06864"
$ws.Range("E19").Value = "  -1.65%  "

$ws.Range("D20").Value = "18.25"
$ws.Range("E20").Value = "  -6.86%  "

$ws.Range("D21").Value = "6.560"
$ws.Range("E21").Value = "  -5.48%  "

$ws.Range("E22").Value = "  -0.22%  "

$ws.Range("D23").Value = "13.01"
$ws.Range("E23").Value = "  -5.12%  "

$ws.Range("D24").Value = "23.135.21"
$ws.Range("E24").Value = "  -3.27%  "

$ws.Range("D25").Value = "2.346"
$ws.Range("E25").Value = "  -4.25%  "

$ws.Range("D26").Value = "2.712"
$ws.Range("E26").Value = "  -7.26%  "

$ws.Range("E27").Value = "  -4.07%  "

$ws.Range("D28").Value = "150.00"
$ws.Range("E28").Value = "  -1.86%  "

$ws.Range("D29").Value = "5.297"
$ws.Range("E29").Value = "  -1.80%  "

$ws.Range("D30").Value = "131.80"
$ws.Range("E30").Value = "  -4.37%  "

$ws.Range("D31").Value = "2.432"
$ws.Range("E31").Value = "  -3.39%  "

$ws.Range("D32").Value = "6.831"
$ws.Range("E32").Value = "  -12.37%  "

$ws.Range("D33").Value = "1.779.11"
$ws.Range("E33").Value = "  -2.74%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "0.9491"
$ws.Range("E34").Value = "  -6.44%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.07672"
$ws.Range("E35").Value = "  -4.40%  "

$ws.Range("D36").Value = "0.02735"
$ws.Range("E36").Value = "  -5.96%  "

$ws.Range("D37").Value = "6.237"
$ws.Range("E37").Value = "  -7.19%  "

$ws.Range("E38").Value = "  -4.40%  "

$ws.Range("D39").Value = "0.08898"
$ws.Range("E39").Value = "  -1.75%  "

$ws.Range("D40").Value = "10.04"
$ws.Range("E40").Value = "  -5.97%  "

$ws.Range("D41").Value = "1.386"
$ws.Range("E41").Value = "  -2.22%  "

$ws.Range("D42").Value = "0.7094"
$ws.Range("E42").Value = "  -6.14%  "

$ws.Range("D43").Value = "12.64"
$ws.Range("E43").Value = "  -5.02%  "

$ws.Range("D44").Value = "15.53"
$ws.Range("E44").Value = "  -4.39%  "

$ws.Range("D45").Value = "0.6605"
$ws.Range("E45").Value = "  -4.88%  "

$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  -0.11%  "

$ws.Range("D47").Value = "2.300"
$ws.Range("E47").Value = "  -5.63%  "

$ws.Range("D48").Value = "3.979"
$ws.Range("E48").Value = "  -2.65%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "131.86"
$ws.Range("E49").Value = "  -2.11%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.07975"
$ws.Range("E50").Value = "  -3.98%  "

$ws.Range("D51").Value = "1.210"
$ws.Range("E51").Value = "  -0.99%  "

# Restore original (unstyled) number format for columns D and E.
$ws.Range("D2:D51").Style = $origStyleD
$ws.Range("E2:E51").Style = $origStyleE
